# Update the header date paragraph
$d = $word.ActiveDocument
$d.Content.Find.Execute('2025-03-31 Monday', $true, $false, $false, $false, $false, $true, 1, $false, '2025-04-01 Tuesday', 2) | Out-Null

# Update each answer cell in the table, in row-major reading order,
# matching the position of each <w:t> run in the document (positional
# replacement is required because several "before" values repeat with
# different "after" values at different positions).
$newValues = @(
    '9+22=31',
    '34+29=63',
    '85-37=48',
    '56-17=39',
    '88-79=9',
    '74+18=92',
    '96-38=58',
    '8+27=35',
    '18+5=23',
    '91-79=12',
    '12-9=3',
    '46+15=61',
    '34+9=43',
    '49+8=57',
    '50-45=5',
    '7+34=41',
    '29+62=91',
    '47+44=91',
    '16+77=93',
    '62-44=18',
    '24-5=19',
    '83-8=75',
    '36-17=19',
    '68+25=93',
    '31-26=5',
    '82-55=27',
    '18+58=76',
    '21-7=14',
    '54-46=8',
    '41-25=16',
    '95-29=66',
    '33-28=5',
    '51-18=33',
    '34-18=16',
    '60-11=49',
    '5+49=54',
    '53-5=48',
    '8+56=64',
    '58+27=85',
    '60-23=37',
    '19+32=51',
    '43-15=28',
    '60-45=15',
    '50-13=37',
    '91-7=84',
    '19+53=72',
    '26+25=51',
    '69+14=83',
    '92-47=45',
    '28-19=9',
    '60-25=35',
    '86-47=39',
    '96-39=57',
    '75-59=16',
    '18+6=24',
    '48-19=29',
    '66+8=74',
    '35-28=7',
    '18+77=95',
    '75+17=92',
    '92-75=17',
    '43-36=7',
    '16-7=9',
    '49+14=63',
    '59+12=71',
    '54-47=7',
    '7+36=43',
    '84-38=46',
    '65+17=82',
    '63-56=7',
    '80-65=15',
    '47+7=54',
    '50-45=5',
    '68+7=75',
    '26+49=75',
    '18+68=86',
    '58+16=74',
    '50-2=48',
    '73-24=49',
    '20-13=7',
    '49+22=71',
    '18+7=25',
    '95-36=59',
    '38+8=46',
    '16+48=64',
    '71-12=59',
    '18+13=31',
    '36+55=91',
    '63-59=4',
    '16+68=84',
    '29+65=94',
    '44+28=72',
    '71-36=35',
    '34+49=83',
    '27+14=41',
    '85-77=8',
    '35-6=29',
    '6+79=85',
    '59+4=63',
    '61-46=15'
)

$t = $d.Tables.Item(1)
$n = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$n]
        $n = $n + 1
    }
}

Write-Output ("Updated " + $n + " cells")
